$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in rows 5 and 6 (columns C and D)
$ws.Range("C5").Value = 0.0199743
$ws.Range("D5").Value = 0.0026535
$ws.Range("C6").Value = 0.0177186
$ws.Range("D6").Value = 0.0048243

# Add new rows 7-10 with values for columns C and D
$ws.Range("C7").Value = 0.0057744
$ws.Range("D7").Value = 0.0062909

$ws.Range("C8").Value = 0.0188013
$ws.Range("D8").Value = 0.0033218

$ws.Range("C9").Value = 0.019903
$ws.Range("D9").Value = 0.0026143

$ws.Range("C10").Value = 0.0199379
$ws.Range("D10").Value = 0.0042961
